$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Team 11 (column D, rows 17-31) - fill in the placeholder roster with the
# real team data, one member "card" (name / e-mail / github id) at a time.

# Team name
$ws.Range("D18").Value = "Git Good"

# Member 1
$ws.Range("D20").Value = "Rafael Oliveira"
$ws.Range("D24").Value = "rcortezb@uci.edu"
$ws.Hyperlinks.Add($ws.Range("D24"), "mailto:rcortezb@uci.edu")
$ws.Range("D28").Value = "rafbel"

# Member 2
$ws.Range("D21").Value = "Hyun Jay Yang"
$ws.Range("D25").Value = "hjyang1@uci.edu"
$ws.Hyperlinks.Add($ws.Range("D25"), "mailto:hjyang1@uci.edu")
$ws.Range("D29").Value = "hjayyang94"

# Member 3
$ws.Range("D22").Value = "Nicolas Grantham"
$ws.Range("D26").Value = "ngrantha@uci.edu"
$ws.Hyperlinks.Add($ws.Range("D26"), "mailto:ngrantha@uci.edu")
$ws.Range("D30").Value = "GranthamAnthem"

# Member 4
$ws.Range("D23").Value = "Chris Zhang"
$ws.Range("D27").Value = "czhang29@uci.edu"
$ws.Hyperlinks.Add($ws.Range("D27"), "mailto:czhang29@uci.edu")
$ws.Range("D31").Value = "ch-zha"

# Leave the selection where the user finished typing
[void]$ws.Range("E26").Select()
